# Commit: "added president's name to all approval excel files"
#
# Adds a new "president" column (F) to the approval-poll sheet and fills
# it with "Trump" for every existing data row, so each row records which
# president the approval/disapproval numbers belong to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "president", styled like the other header cells
# (A1:E1) so it keeps the same bold / centered look.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "president"

# F2:F101 = "Trump" for every existing data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 101 }
$ws.Range("F2:F$lastRow").Value = "Trump"
